$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Range("G16:G20").NumberFormat = "@"

$ws.Range("E16").Value = 1
$ws.Range("G16").Value = "3.23 %"

$ws.Range("E17").Value = 27
$ws.Range("G17").Value = "87.1 %"

$ws.Range("E19").Value = 2
$ws.Range("G19").Value = "6.45 %"

$ws.Range("E20").Value = 1
$ws.Range("G20").Value = "3.23 %"

$ws.Range("G16:G20").NumberFormat = "General"
